$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 82
$ws.Range("H82").Value = 741.7143
$ws.Range("I82").Value = 741.7143
$ws.Range("K82").Value = 2225.1429
$ws.Range("M82").Value = -1819.1429
# row 85
$ws.Range("H85").Value = 741.7143
$ws.Range("I85").Value = 741.7143
$ws.Range("K85").Value = 2225.1429
$ws.Range("M85").Value = -821.1428999999998
# row 112
$ws.Range("H112").Value = 711172.4
$ws.Range("J112").Value = 784704.0600000001
$ws.Range("L112").Value = 2354112.18
$ws.Range("N112").Value = -2356328.18
# row 129
$ws.Range("H129").Value = 880.2
$ws.Range("J129").Value = 1074.9524
$ws.Range("L129").Value = 3224.857199999999
$ws.Range("N129").Value = -13224.8572

$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 1342.6
$ws.Range("I2").Value = 1300
$ws.Range("K2").Value = 1300
$ws.Range("M2").Value = -1187
# row 46
$ws.Range("H46").Value = 1252.8334
$ws.Range("I46").Value = 1137
$ws.Range("J46").Value = 1276
$ws.Range("K46").Value = 1137
$ws.Range("L46").Value = 1276
$ws.Range("M46").Value = -818
$ws.Range("N46").Value = -1914
# row 74
$ws.Range("H74").Value = 280274.75
$ws.Range("I74").Value = 2617.4285
$ws.Range("J74").Value = 1252075.4
$ws.Range("K74").Value = 2617.4285
$ws.Range("L74").Value = 1252075.4
$ws.Range("M74").Value = -1743.4285
$ws.Range("N74").Value = -1253823.4
# row 77
$ws.Range("H77").Value = 280274.75
$ws.Range("I77").Value = 2617.4285
$ws.Range("J77").Value = 1252075.4
$ws.Range("K77").Value = 13087.1425
$ws.Range("L77").Value = 6260377
$ws.Range("M77").Value = -8719.1425
$ws.Range("N77").Value = -6269113
# row 98
$ws.Range("H98").Value = 17699.75
$ws.Range("J98").Value = 17699.75
$ws.Range("L98").Value = 17699.75
$ws.Range("N98").Value = -23689.75
# row 116
$ws.Range("H116").Value = 1342.6
$ws.Range("I116").Value = 1300
$ws.Range("K116").Value = 1300
$ws.Range("M116").Value = 994
# row 122
$ws.Range("H122").Value = 2747.8809
$ws.Range("J122").Value = 3686.4
$ws.Range("L122").Value = 11059.2
$ws.Range("N122").Value = -15959.2
# row 132
$ws.Range("H132").Value = 19254.203
$ws.Range("I132").Value = 25589.744
$ws.Range("J132").Value = 2227.4375
$ws.Range("K132").Value = 76769.23199999999
$ws.Range("L132").Value = 6682.3125
$ws.Range("M132").Value = -74239.23199999999
$ws.Range("N132").Value = -11742.3125

$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 1342.6
$ws.Range("I3").Value = 1300
$ws.Range("K3").Value = 1300
$ws.Range("M3").Value = -1186
# row 105
$ws.Range("H105").Value = 2340.4375
$ws.Range("I105").Value = 1984.7
$ws.Range("J105").Value = 2933.3333
$ws.Range("K105").Value = 1984.7
$ws.Range("L105").Value = 2933.3333
$ws.Range("M105").Value = -237.7
$ws.Range("N105").Value = -6427.3333

$ws = $wb.Worksheets.Item("CRP")
# row 99
$ws.Range("H99").Value = 7941.2
$ws.Range("I99").Value = 12328
$ws.Range("J99").Value = 5016.6665
$ws.Range("K99").Value = 12328
$ws.Range("L99").Value = 5016.6665
$ws.Range("M99").Value = -10830
$ws.Range("N99").Value = -8012.6665
# row 126
$ws.Range("H126").Value = 7941.2
$ws.Range("I126").Value = 12328
$ws.Range("J126").Value = 5016.6665
$ws.Range("K126").Value = 36984
$ws.Range("L126").Value = 15049.9995
$ws.Range("M126").Value = -34514
$ws.Range("N126").Value = -19989.9995

$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 857.17145
$ws.Range("J5").Value = 1311
$ws.Range("L5").Value = 3933
$ws.Range("N5").Value = -4157
# row 40
$ws.Range("H40").Value = 4787.048
$ws.Range("J40").Value = 14257
$ws.Range("L40").Value = 57028
$ws.Range("N40").Value = -57166
# row 92
$ws.Range("H92").Value = 576.375
$ws.Range("I92").Value = 490.25
$ws.Range("J92").Value = 662.5
$ws.Range("K92").Value = 1470.75
$ws.Range("L92").Value = 1987.5
$ws.Range("M92").Value = -222.75
$ws.Range("N92").Value = -4483.5
# row 98
$ws.Range("H98").Value = 598.8
$ws.Range("I98").Value = 550.75
$ws.Range("J98").Value = 630.8333
$ws.Range("K98").Value = 1652.25
$ws.Range("L98").Value = 1892.4999
$ws.Range("M98").Value = -154.25
$ws.Range("N98").Value = -4888.4999
# row 102
$ws.Range("H102").Value = 4950
$ws.Range("J102").Value = 4933.3335
$ws.Range("L102").Value = 14800.0005
$ws.Range("N102").Value = -19668.0005
# row 107
$ws.Range("H107").Value = 110.23077
$ws.Range("I107").Value = 111.083336
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 333.250008
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = 1586.749992
$ws.Range("N107").Value = -4140
# row 114
$ws.Range("H114").Value = 2241.8823
$ws.Range("I114").Value = 1608.1
$ws.Range("J114").Value = 3147.2856
$ws.Range("K114").Value = 4824.299999999999
$ws.Range("L114").Value = 9441.856800000001
$ws.Range("M114").Value = -1570.299999999999
$ws.Range("N114").Value = -15949.8568
# row 121
$ws.Range("H121").Value = 1383.1025
$ws.Range("I121").Value = 389.14285
$ws.Range("J121").Value = 1939.72
$ws.Range("K121").Value = 1167.42855
$ws.Range("L121").Value = 5819.16
$ws.Range("M121").Value = 142.5714499999999
$ws.Range("N121").Value = -8439.16
# row 123
$ws.Range("H123").Value = 2847.5112
$ws.Range("I123").Value = 1401.4286
$ws.Range("J123").Value = 3113.8948
$ws.Range("K123").Value = 4204.2858
$ws.Range("L123").Value = 9341.6844
$ws.Range("M123").Value = -1754.2858
$ws.Range("N123").Value = -14241.6844
# row 131
$ws.Range("H131").Value = 797.96295
$ws.Range("I131").Value = 456.8
$ws.Range("J131").Value = 998.64703
$ws.Range("K131").Value = 1370.4
$ws.Range("L131").Value = 2995.94109
$ws.Range("M131").Value = 3669.6
$ws.Range("N131").Value = -13075.94109
# row 135
$ws.Range("H135").Value = 857.17145
$ws.Range("J135").Value = 1311
$ws.Range("L135").Value = 11799
$ws.Range("N135").Value = -16869

$ws = $wb.Worksheets.Item("GSM")
# row 39
$ws.Range("H39").Value = 30355.5
$ws.Range("J39").Value = 30355.5
$ws.Range("L39").Value = 30355.5
$ws.Range("N39").Value = -31419.5
# row 122
$ws.Range("H122").Value = 1686.1
$ws.Range("I122").Value = 2018.6154
$ws.Range("J122").Value = 1068.5714
$ws.Range("K122").Value = 6055.8462
$ws.Range("L122").Value = 3205.7142
$ws.Range("M122").Value = -3605.8462
$ws.Range("N122").Value = -8105.7142

$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 2279.0908
$ws.Range("I7").Value = 1685.5555
$ws.Range("K7").Value = 1685.5555
$ws.Range("M7").Value = -1573.5555
# row 22
$ws.Range("H22").Value = 1852582.4
$ws.Range("I22").Value = 2564471.8
$ws.Range("J22").Value = 1670
$ws.Range("K22").Value = 2564471.8
$ws.Range("L22").Value = 1670
$ws.Range("M22").Value = -2564176.8
$ws.Range("N22").Value = -2260
# row 27
$ws.Range("H27").Value = 1852582.4
$ws.Range("I27").Value = 2564471.8
$ws.Range("J27").Value = 1670
$ws.Range("K27").Value = 2564471.8
$ws.Range("L27").Value = 1670
$ws.Range("M27").Value = -2564364.8
$ws.Range("N27").Value = -1884
# row 55
$ws.Range("H55").Value = 243.8158
$ws.Range("I55").Value = 284.33334
$ws.Range("J55").Value = 193.76471
$ws.Range("K55").Value = 284.33334
$ws.Range("L55").Value = 193.76471
$ws.Range("M55").Value = -111.33334
$ws.Range("N55").Value = -539.76471
# row 126
$ws.Range("H126").Value = 2279.0908
$ws.Range("I126").Value = 1685.5555
$ws.Range("K126").Value = 5056.666499999999
$ws.Range("M126").Value = -2586.666499999999
# row 136
$ws.Range("H136").Value = 6906.1924
$ws.Range("I136").Value = 8275.706
$ws.Range("J136").Value = 4319.3335
$ws.Range("K136").Value = 24827.118
$ws.Range("L136").Value = 12958.0005
$ws.Range("M136").Value = -22277.118
$ws.Range("N136").Value = -18058.0005

$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 1467.5588
$ws.Range("I132").Value = 827.5599999999999
$ws.Range("J132").Value = 3245.3333
$ws.Range("K132").Value = 2482.68
$ws.Range("L132").Value = 9735.999899999999
$ws.Range("M132").Value = 47.32000000000016
$ws.Range("N132").Value = -14795.9999
